# Automatische test-sync: 2025-08-05 16:46:50
# Append a new test-mail log row to "Logs" and refresh the "Dashboard"
# category summary table that feeds the bar chart.

$wb = $excel.ActiveWorkbook

$logs = $wb.Worksheets.Item("Logs")
$dashboard = $wb.Worksheets.Item("Dashboard")

# --- Logs: append row 8 ------------------------------------------------
$newRow = 8

$logs.Cells.Item($newRow, 1).Value = "Kun jij dit even regelen?"
$logs.Cells.Item($newRow, 2).Value = "mailmind.test@zohomail.eu"
$logs.Cells.Item($newRow, 3).Value = "Testmail #1: Kun jij dit even regelen?"
$logs.Cells.Item($newRow, 4).Value = "Planning / Afspraak"
$logs.Cells.Item($newRow, 5).Value = "Bedankt, we hebben dit doorgestuurd naar planning@bedrijf.nl."
$logs.Cells.Item($newRow, 6).Value = "2025-08-05 16:46:17"
$logs.Cells.Item($newRow, 7).Value = "Ja"
$logs.Cells.Item($newRow, 8).Value = "Ja"
$logs.Cells.Item($newRow, 9).Value = "Nee"
$logs.Cells.Item($newRow, 10).Value = "Nee"

# --- Dashboard: refresh category counts (rows 3-5 reordered) -----------
$dashboard.Cells.Item(3, 1).Value = "Planning / Afspraak"
$dashboard.Cells.Item(3, 2).Value = 2

$dashboard.Cells.Item(4, 1).Value = "Retour / Terugbetaling"
$dashboard.Cells.Item(4, 2).Value = 1

$dashboard.Cells.Item(5, 1).Value = "Klacht / Probleem"
$dashboard.Cells.Item(5, 2).Value = 1

# --- Logs: extend the conditional-formatting ranges to cover row 8 -----
$logs.Range("D2:D7").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("D2:D8"))
$logs.Range("G2:G7").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("G2:G8"))
$logs.Range("H2:H7").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("H2:H8"))
$logs.Range("I2:I7").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("I2:I8"))
$logs.Range("J2:J7").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("J2:J8"))
